# Update Sheets via scheduled runner: refresh leve crafting profit data
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) for affected leves.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 366.42856
$ws.Range("I12").Value = 260.83334
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 260.83334
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -90.83334000000002
$ws.Range("N12").Value = -1340
# Row 17
$ws.Range("H17").Value = 2074.5557
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2074.5557
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6223.6671
$ws.Range("N17").Value = -6559.6671
# Row 19
$ws.Range("H19").Value = 3600.5
$ws.Range("I19").Value = 3601
$ws.Range("J19").Value = 3600
$ws.Range("K19").Value = 3601
$ws.Range("L19").Value = 3600
$ws.Range("M19").Value = -3426
$ws.Range("N19").Value = -3950
# Row 28
$ws.Range("H28").Value = 1326.6666
$ws.Range("I28").Value = 993
$ws.Range("J28").Value = 2995
$ws.Range("K28").Value = 993
$ws.Range("L28").Value = 2995
$ws.Range("M28").Value = -508
# Row 62
$ws.Range("H62").Value = 283
$ws.Range("I62").Value = 283
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 283
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 341
# Row 65
$ws.Range("H65").Value = 283
$ws.Range("I65").Value = 283
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1415
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 1705
# Row 82
$ws.Range("H82").Value = 500
$ws.Range("I82").Value = 500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1094
# Row 85
$ws.Range("H85").Value = 500
$ws.Range("I85").Value = 500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -96
# Row 116
$ws.Range("H116").Value = 4273
$ws.Range("I116").Value = 2300
$ws.Range("J116").Value = 4667.6
$ws.Range("K116").Value = 2300
$ws.Range("L116").Value = 4667.6
$ws.Range("M116").Value = 1142
$ws.Range("N116").Value = -11551.6
# Row 138
$ws.Range("H138").Value = 4352.6904
$ws.Range("I138").Value = 2558.182
$ws.Range("J138").Value = 4989.4517
$ws.Range("K138").Value = 7674.545999999999
$ws.Range("L138").Value = 14968.3551
$ws.Range("M138").Value = -2534.545999999999
$ws.Range("N138").Value = -25248.3551
# Row 141
$ws.Range("H141").Value = 2975.1667
$ws.Range("I141").Value = 2581.3809
$ws.Range("J141").Value = 5731.6665
$ws.Range("K141").Value = 7744.1427
$ws.Range("L141").Value = 17194.9995
$ws.Range("M141").Value = -2564.1427

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4746.079
$ws.Range("I32").Value = 4675.0586
$ws.Range("J32").Value = 5349.75
$ws.Range("K32").Value = 4675.0586
$ws.Range("L32").Value = 5349.75
$ws.Range("M32").Value = -4388.0586
$ws.Range("N32").Value = -5923.75
# Row 50
$ws.Range("H50").Value = 11485.5
$ws.Range("I50").Value = 10296.333
$ws.Range("J50").Value = 15053
$ws.Range("K50").Value = 10296.333
$ws.Range("L50").Value = 15053
$ws.Range("M50").Value = -9582.333000000001
$ws.Range("N50").Value = -16481
# Row 122
$ws.Range("H122").Value = 2898.5
$ws.Range("I122").Value = 2898.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8695.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6245.5
# Row 132
$ws.Range("H132").Value = 3608.0476
$ws.Range("I132").Value = 2616
$ws.Range("J132").Value = 4699.3
$ws.Range("K132").Value = 7848
$ws.Range("L132").Value = 14097.9
$ws.Range("M132").Value = -5318

$ws = $wb.Worksheets.Item("BSM")
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 131
$ws.Range("H131").Value = 27496
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 27496
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 27496
$ws.Range("N131").Value = -37576
# Row 134
$ws.Range("H134").Value = 3021.158
$ws.Range("I134").Value = 3106.1667
$ws.Range("J134").Value = 1491
$ws.Range("K134").Value = 9318.500100000001
$ws.Range("L134").Value = 4473
$ws.Range("M134").Value = -6783.500100000001
$ws.Range("N134").Value = -9543

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 579.1667
$ws.Range("I22").Value = 798.3333
$ws.Range("J22").Value = 506.1111
$ws.Range("K22").Value = 798.3333
$ws.Range("L22").Value = 506.1111
$ws.Range("M22").Value = -448.3333
$ws.Range("N22").Value = -1206.1111
# Row 58
$ws.Range("H58").Value = 2382.6
$ws.Range("I58").Value = 1961.6666
$ws.Range("J58").Value = 3014
$ws.Range("K58").Value = 1961.6666
$ws.Range("L58").Value = 3014
$ws.Range("M58").Value = -1758.6666
# Row 99
$ws.Range("H99").Value = 4750
$ws.Range("I99").Value = 4750
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4750
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3252
$ws.Range("N99").ClearContents()
# Row 105
$ws.Range("H105").Value = 4259.2
$ws.Range("I105").Value = 3065.3333
$ws.Range("J105").Value = 6050
$ws.Range("K105").Value = 3065.3333
$ws.Range("L105").Value = 6050
$ws.Range("M105").Value = -1318.3333
# Row 126
$ws.Range("H126").Value = 4750
$ws.Range("I126").Value = 4750
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14250
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11780
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 3149.2104
$ws.Range("I132").Value = 2780.5715
$ws.Range("J132").Value = 4181.4
$ws.Range("K132").Value = 8341.7145
$ws.Range("L132").Value = 12544.2
$ws.Range("M132").Value = -5811.7145
$ws.Range("N132").Value = -17604.2
# Row 134
$ws.Range("H134").Value = 1664.5
$ws.Range("I134").Value = 1866.9166
$ws.Range("J134").Value = 450
$ws.Range("K134").Value = 5600.7498
$ws.Range("L134").Value = 1350
$ws.Range("M134").Value = -3065.7498
$ws.Range("N134").Value = -6420
# Row 136
$ws.Range("H136").Value = 2382.6
$ws.Range("I136").Value = 1961.6666
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 5884.9998
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -3334.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 86.666664
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 360
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -247
$ws.Range("N2").Value = -826
# Row 38
$ws.Range("H38").Value = 250
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 750
$ws.Range("N38").Value = -1444
$ws.Range("M38").ClearContents()
# Row 40
$ws.Range("H40").Value = 209.85715
$ws.Range("I40").Value = 161.5
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 646
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -577
# Row 86
$ws.Range("H86").Value = 296.33334
$ws.Range("I86").Value = 296.33334
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 889.0000200000001
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 296.9999799999999
# Row 89
$ws.Range("H89").Value = 296.33334
$ws.Range("I89").Value = 296.33334
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2667.00006
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3260.99994
# Row 92
$ws.Range("H92").Value = 297.5
$ws.Range("I92").Value = 295
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 885
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 363
$ws.Range("N92").Value = -3396
# Row 113
$ws.Range("H113").Value = 1012.6667
$ws.Range("I113").Value = 1224.4
$ws.Range("J113").Value = 946.5
$ws.Range("K113").Value = 3673.2
$ws.Range("L113").Value = 2839.5
$ws.Range("M113").Value = -1503.2
$ws.Range("N113").Value = -7179.5
# Row 133
$ws.Range("H133").Value = 7843.3335
$ws.Range("I133").Value = 7843.3335
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 23530.0005
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -18470.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1219.4
$ws.Range("I122").Value = 1499.5
$ws.Range("J122").Value = 1032.6666
$ws.Range("K122").Value = 4498.5
$ws.Range("L122").Value = 3097.9998
$ws.Range("M122").Value = -2048.5
$ws.Range("N122").Value = -7997.9998
# Row 126
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -19940
# Row 132
$ws.Range("H132").Value = 1228.871
$ws.Range("I132").Value = 657.6539
$ws.Range("J132").Value = 4199.2
$ws.Range("K132").Value = 1972.9617
$ws.Range("L132").Value = 12597.6
$ws.Range("M132").Value = 557.0382999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3169.8
$ws.Range("I40").Value = 3349.75
$ws.Range("J40").Value = 2450
$ws.Range("K40").Value = 3349.75
$ws.Range("L40").Value = 2450
$ws.Range("M40").Value = -3213.75
$ws.Range("N40").Value = -2722
# Row 46
$ws.Range("H46").Value = 240
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 240
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 240
$ws.Range("N46").Value = -616
$ws.Range("M46").ClearContents()
# Row 61
$ws.Range("H61").Value = 2089.4285
$ws.Range("I61").Value = 1771
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1771
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1569
# Row 68
$ws.Range("H68").Value = 2925
$ws.Range("I68").Value = 2925
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2925
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2176
# Row 71
$ws.Range("H71").Value = 2925
$ws.Range("I71").Value = 2925
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14625
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -10881
# Row 113
$ws.Range("H113").Value = 2089.4285
$ws.Range("I113").Value = 1771
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1771
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 399
# Row 122
$ws.Range("H122").Value = 3002.7144
$ws.Range("I122").Value = 3002.8333
$ws.Range("J122").Value = 3002
$ws.Range("K122").Value = 9008.499899999999
$ws.Range("L122").Value = 9006
$ws.Range("M122").Value = -6558.499899999999
$ws.Range("N122").Value = -13906
# Row 132
$ws.Range("H132").Value = 3215.75
$ws.Range("I132").Value = 2476
$ws.Range("J132").Value = 4090
$ws.Range("K132").Value = 7428
$ws.Range("L132").Value = 12270
$ws.Range("M132").Value = -4898
$ws.Range("N132").Value = -17330

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1638.8572
$ws.Range("I122").Value = 1638.8572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4916.571599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2466.571599999999
# Row 132
$ws.Range("H132").Value = 1909.0702
$ws.Range("I132").Value = 1238.9286
$ws.Range("J132").Value = 3785.4666
$ws.Range("K132").Value = 3716.7858
$ws.Range("L132").Value = 11356.3998
$ws.Range("M132").Value = -1186.7858
# Row 136
$ws.Range("H136").Value = 1392.2778
$ws.Range("I136").Value = 1042.9166
$ws.Range("J136").Value = 2091
$ws.Range("K136").Value = 3128.7498
$ws.Range("L136").Value = 6273
$ws.Range("M136").Value = -578.7498000000001
